$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 25007498
$ws.Range("I111").Value = 31257632
$ws.Range("K111").Value = 93772896
$ws.Range("M111").Value = -93769829

$ws.Range("H113").Value = 45463760
$ws.Range("I113").Value = 3582.4
$ws.Range("J113").Value = 83347240
$ws.Range("K113").Value = 3582.4
$ws.Range("L113").Value = 83347240
$ws.Range("M113").Value = -328.4000000000001
$ws.Range("N113").Value = -83353748

$ws.Range("H132").Value = 1148.079
$ws.Range("I132").Value = 1201.0322
$ws.Range("K132").Value = 3603.0966
$ws.Range("M132").Value = -1073.0966

$ws.Range("H138").Value = 5588.237
$ws.Range("I138").Value = 921.25
$ws.Range("J138").Value = 10773.777
$ws.Range("K138").Value = 2763.75
$ws.Range("L138").Value = 32321.331
$ws.Range("M138").Value = 2376.25
$ws.Range("N138").Value = -42601.331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3128733
$ws.Range("I32").Value = 3381901.5
$ws.Range("K32").Value = 3381901.5
$ws.Range("M32").Value = -3381614.5

$ws.Range("H45").Value = 2980.9443
$ws.Range("I45").Value = 1699.7778
$ws.Range("J45").Value = 4262.1113
$ws.Range("K45").Value = 1699.7778
$ws.Range("L45").Value = 4262.1113
$ws.Range("M45").Value = -1322.7778
$ws.Range("N45").Value = -5016.1113

$ws.Range("H61").Value = 6125.1763
$ws.Range("I61").Value = 3401.7273
$ws.Range("J61").Value = 11118.167
$ws.Range("K61").Value = 3401.7273
$ws.Range("L61").Value = 11118.167
$ws.Range("M61").Value = -3189.7273
$ws.Range("N61").Value = -11542.167

$ws.Range("H74").Value = 35890.105
$ws.Range("I74").Value = 45681.027
$ws.Range("J74").Value = 3847.0908
$ws.Range("K74").Value = 45681.027
$ws.Range("L74").Value = 3847.0908
$ws.Range("M74").Value = -44807.027
$ws.Range("N74").Value = -5595.0908

$ws.Range("H77").Value = 35890.105
$ws.Range("I77").Value = 45681.027
$ws.Range("J77").Value = 3847.0908
$ws.Range("K77").Value = 228405.135
$ws.Range("L77").Value = 19235.454
$ws.Range("M77").Value = -224037.135
$ws.Range("N77").Value = -27971.454

$ws.Range("H88").Value = 2318.353
$ws.Range("J88").Value = 2368
$ws.Range("L88").Value = 2368
$ws.Range("N88").Value = -3180

$ws.Range("H91").Value = 2318.353
$ws.Range("J91").Value = 2368
$ws.Range("L91").Value = 2368
$ws.Range("N91").Value = -5176

$ws.Range("H122").Value = 13184.863
$ws.Range("I122").Value = 13845.632
$ws.Range("K122").Value = 41536.896
$ws.Range("M122").Value = -39086.896

$ws.Range("H132").Value = 5451.525
$ws.Range("I132").Value = 5443.2915
$ws.Range("K132").Value = 16329.8745
$ws.Range("M132").Value = -13799.8745

$ws.Range("H136").Value = 6125.1763
$ws.Range("I136").Value = 3401.7273
$ws.Range("J136").Value = 11118.167
$ws.Range("K136").Value = 10205.1819
$ws.Range("L136").Value = 33354.501
$ws.Range("M136").Value = -7655.1819
$ws.Range("N136").Value = -38454.501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 62504404
$ws.Range("I86").Value = 3540.9167
$ws.Range("K86").Value = 3540.9167
$ws.Range("M86").Value = -2417.9167

$ws.Range("H89").Value = 62504404
$ws.Range("I89").Value = 3540.9167
$ws.Range("K89").Value = 17704.5835
$ws.Range("M89").Value = -12088.5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 56353.332
$ws.Range("J18").Value = 56353.332
$ws.Range("L18").Value = 56353.332
$ws.Range("N18").Value = -56813.332

$ws.Range("H31").Value = 6202.143
$ws.Range("I31").Value = 3690.3462
$ws.Range("J31").Value = 8379.032999999999
$ws.Range("K31").Value = 3690.3462
$ws.Range("L31").Value = 8379.032999999999
$ws.Range("M31").Value = -3395.3462
$ws.Range("N31").Value = -8969.032999999999

$ws.Range("H34").Value = 6202.143
$ws.Range("I34").Value = 3690.3462
$ws.Range("J34").Value = 8379.032999999999
$ws.Range("K34").Value = 3690.3462
$ws.Range("L34").Value = 8379.032999999999
$ws.Range("M34").Value = -3488.3462
$ws.Range("N34").Value = -8783.032999999999

$ws.Range("H82").Value = 56181
$ws.Range("J82").Value = 56181
$ws.Range("L82").Value = 56181
$ws.Range("N82").Value = -56903

$ws.Range("H85").Value = 56181
$ws.Range("J85").Value = 56181
$ws.Range("L85").Value = 56181
$ws.Range("N85").Value = -58677

$ws.Range("H117").Value = 59999.5
$ws.Range("J117").Value = 59999.5
$ws.Range("L117").Value = 59999.5
$ws.Range("N117").Value = -69177.5

$ws.Range("H122").Value = 1370.75
$ws.Range("I122").Value = 1370.75
$ws.Range("K122").Value = 4112.25
$ws.Range("M122").Value = -1662.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 1000000000
$ws.Range("I76").Value = 1000000000
$ws.Range("K76").Value = 3000000000
$ws.Range("M76").Value = -2999999617

$ws.Range("H79").Value = 1000000000
$ws.Range("I79").Value = 1000000000
$ws.Range("K79").Value = 3000000000
$ws.Range("M79").Value = -2999998674

$ws.Range("H107").Value = 15385172
$ws.Range("J107").Value = 20000558
$ws.Range("L107").Value = 60001674
$ws.Range("N107").Value = -60005514

$ws.Range("H132").Value = 21998
$ws.Range("J132").Value = 24247.5
$ws.Range("L132").Value = 218227.5
$ws.Range("N132").Value = -223287.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7234.1514
$ws.Range("I70").Value = 5705.65
$ws.Range("K70").Value = 5705.65
$ws.Range("M70").Value = -5435.65

$ws.Range("H73").Value = 7234.1514
$ws.Range("I73").Value = 5705.65
$ws.Range("K73").Value = 5705.65
$ws.Range("M73").Value = -4769.65

$ws.Range("H97").Value = 2792.5217
$ws.Range("J97").Value = 2939.9
$ws.Range("L97").Value = 2939.9
$ws.Range("N97").Value = -3931.9

$ws.Range("H101").Value = 49964.332
$ws.Range("J101").Value = 49964.332
$ws.Range("L101").Value = 49964.332
$ws.Range("N101").Value = -56454.332

$ws.Range("H122").Value = 41110.68
$ws.Range("I122").Value = 75941.42999999999
$ws.Range("K122").Value = 227824.29
$ws.Range("M122").Value = -225374.29

$ws.Range("H123").Value = 43334
$ws.Range("J123").Value = 43334
$ws.Range("L123").Value = 43334
$ws.Range("N123").Value = -48234

$ws.Range("H126").Value = 3089.8
$ws.Range("I126").Value = 3071.4
$ws.Range("J126").Value = 3099
$ws.Range("K126").Value = 9214.200000000001
$ws.Range("L126").Value = 9297
$ws.Range("M126").Value = -6744.200000000001
$ws.Range("N126").Value = -14237

$ws.Range("H132").Value = 1924.6666
$ws.Range("I132").Value = 1828.0975
$ws.Range("J132").Value = 2490.2856
$ws.Range("K132").Value = 5484.2925
$ws.Range("L132").Value = 7470.8568
$ws.Range("M132").Value = -2954.2925
$ws.Range("N132").Value = -12530.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5317.6113
$ws.Range("I40").Value = 4534.4443
$ws.Range("K40").Value = 4534.4443
$ws.Range("M40").Value = -4398.4443

$ws.Range("H46").Value = 13891018
$ws.Range("J46").Value = 13891018
$ws.Range("L46").Value = 13891018
$ws.Range("N46").Value = -13891394

$ws.Range("H93").Value = 2987.4285
$ws.Range("I93").Value = 2885.2856
$ws.Range("J93").Value = 3140.6428
$ws.Range("K93").Value = 2885.2856
$ws.Range("L93").Value = 3140.6428
$ws.Range("M93").Value = -1637.2856
$ws.Range("N93").Value = -5636.6428

$ws.Range("H122").Value = 3948.913
$ws.Range("I122").Value = 2937.5
$ws.Range("J122").Value = 5522.222
$ws.Range("K122").Value = 8812.5
$ws.Range("L122").Value = 16566.666
$ws.Range("M122").Value = -6362.5
$ws.Range("N122").Value = -21466.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 10000
$ws.Range("I52").Value = 10000
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("M52").Value = -9774

$ws.Range("H74").Value = 2111
$ws.Range("J74").Value = 2111
$ws.Range("L74").Value = 2111
$ws.Range("N74").Value = -3983

$ws.Range("H77").Value = 2111
$ws.Range("J77").Value = 2111
$ws.Range("L77").Value = 6333
$ws.Range("N77").Value = -15693

$ws.Range("H107").Value = 1108.7778
$ws.Range("J107").Value = 1375
$ws.Range("L107").Value = 4125
$ws.Range("N107").Value = -7965

$ws.Range("H122").Value = 4382.1816
$ws.Range("I122").Value = 2315
$ws.Range("K122").Value = 6945
$ws.Range("M122").Value = -4495

$ws.Range("H126").Value = 3728.4
$ws.Range("I126").Value = 2175.0908
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 6525.2724
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -4055.2724
$ws.Range("N126").Value = -28940
